# Update viewer/attendance counts (column F) across the three sheets that
# carry event data: 展览 (Exhibitions), 演出 (Performances) and
# 全部类型 (All types - combined view). 本地生活 (Local life) has no data
# rows and needs no changes.

$wb = $excel.ActiveWorkbook

# ---- 展览 (sheet1) ----
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 82
$ws.Range("F5").Value = 281
$ws.Range("F7").Value = 103
$ws.Range("F9").Value = 22
$ws.Range("F10").Value = 52
$ws.Range("F12").Value = 118
$ws.Range("F13").Value = 2479
$ws.Range("F14").Value = 55
$ws.Range("F15").Value = 26
$ws.Range("F16").Value = 61
$ws.Range("F17").Value = 14
$ws.Range("F18").Value = 42
$ws.Range("F19").Value = 537
$ws.Range("F20").Value = 585
$ws.Range("F22").Value = 90
$ws.Range("F25").Value = 2071
$ws.Range("F26").Value = 4174
$ws.Range("F29").Value = 467
$ws.Range("F30").Value = 1221
$ws.Range("F31").Value = 239
$ws.Range("F32").Value = 2125
$ws.Range("F33").Value = 567
$ws.Range("F34").Value = 473
$ws.Range("F35").Value = 67
$ws.Range("F36").Value = 125
$ws.Range("F38").Value = 435
$ws.Range("F39").Value = 720
$ws.Range("F42").Value = 5
$ws.Range("F43").Value = 431

# ---- 演出 (sheet2) ----
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 46

# ---- 全部类型 (sheet4) ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 82
$ws.Range("F5").Value = 281
$ws.Range("F7").Value = 103
$ws.Range("F9").Value = 22
$ws.Range("F10").Value = 52
$ws.Range("F12").Value = 118
$ws.Range("F13").Value = 2479
$ws.Range("F14").Value = 55
$ws.Range("F15").Value = 26
$ws.Range("F16").Value = 61
$ws.Range("F17").Value = 46
$ws.Range("F18").Value = 14
$ws.Range("F19").Value = 42
$ws.Range("F20").Value = 537
$ws.Range("F21").Value = 585
$ws.Range("F23").Value = 90
$ws.Range("F26").Value = 2071
$ws.Range("F27").Value = 4174
$ws.Range("F30").Value = 467
$ws.Range("F31").Value = 1221
$ws.Range("F32").Value = 239
$ws.Range("F33").Value = 2125
$ws.Range("F34").Value = 567
$ws.Range("F35").Value = 473
$ws.Range("F36").Value = 67
$ws.Range("F37").Value = 125
$ws.Range("F39").Value = 435
$ws.Range("F40").Value = 720
$ws.Range("F43").Value = 5
$ws.Range("F44").Value = 431

$wb.Save()
